$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-07-24 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-25 Thursday", 2)

# Update the table of division problems (Table 1, rows 1,5,9,13,17 have content)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "13÷2=6, 1"
$t.Cell(1,2).Range.Text  = "92÷2=46, 0"
$t.Cell(1,3).Range.Text  = "40÷2=20, 0"
$t.Cell(1,4).Range.Text  = "83÷8=10, 3"
$t.Cell(1,5).Range.Text  = "55÷4=13, 3"

$t.Cell(5,1).Range.Text  = "89÷9=9, 8"
$t.Cell(5,2).Range.Text  = "73÷7=10, 3"
$t.Cell(5,3).Range.Text  = "81÷3=27, 0"
$t.Cell(5,4).Range.Text  = "94÷9=10, 4"
$t.Cell(5,5).Range.Text  = "23÷2=11, 1"

$t.Cell(9,1).Range.Text  = "76÷2=38, 0"
$t.Cell(9,2).Range.Text  = "41÷3=13, 2"
$t.Cell(9,3).Range.Text  = "73÷7=10, 3"
$t.Cell(9,4).Range.Text  = "72÷8=9, 0"
$t.Cell(9,5).Range.Text  = "56÷6=9, 2"

$t.Cell(13,1).Range.Text = "60÷5=12, 0"
$t.Cell(13,2).Range.Text = "98÷2=49, 0"
$t.Cell(13,3).Range.Text = "33÷3=11, 0"
$t.Cell(13,4).Range.Text = "91÷6=15, 1"
$t.Cell(13,5).Range.Text = "68÷8=8, 4"

$t.Cell(17,1).Range.Text = "52÷3=17, 1"
$t.Cell(17,2).Range.Text = "49÷5=9, 4"
$t.Cell(17,3).Range.Text = "33÷3=11, 0"
$t.Cell(17,4).Range.Text = "79÷8=9, 7"
$t.Cell(17,5).Range.Text = "16÷2=8, 0"
